# Generate Report for Handback
# Adds a new handback record (d8893c56-fe19-464e-94ef-f52a94eba1ad.md) as
# row 5 of all three sheets: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

$fileName      = "d8893c56-fe19-464e-94ef-f52a94eba1ad.md"
$pathAndName   = "e2e\d8893c56-fe19-464e-94ef-f52a94eba1ad.md"
$ext           = ".md"
$statusInSync  = "Handed back: in sync with en-US"

$zhXlf   = "d8893c56-fe19-464e-94ef-f52a94eba1ad.ce58ac674dba2c9d903ff8d34ca36f55a512b5c0.zh-cn.xlf"
$deXlf   = "d8893c56-fe19-464e-94ef-f52a94eba1ad.ce58ac674dba2c9d903ff8d34ca36f55a512b5c0.de-de.xlf"

$overviewDate = "2017-02-21 11:08:44"

$zhHandoffDate   = "2017-02-21 11:08:27"
$zhHandbackDate  = "2017-02-21 11:09:26"

$deHandoffDate   = "2017-02-21 11:08:44"
$deHandbackDate  = "2017-02-21 11:09:48"

$dateFmt = "yyyy-mm-dd HH:mm:ss"

# -------------------------------------------------------------------
# Sheet 1: Overview
# -------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item(1)
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A5").Value = $fileName
$wsOverview.Range("B5").Value = $pathAndName
$wsOverview.Range("C5").Value = $ext
$wsOverview.Range("E5").Value = $statusInSync
$wsOverview.Range("F5").Value = $statusInSync
$wsOverview.Range("G5").Value = $overviewDate
$wsOverview.Range("G5").NumberFormat = $dateFmt

$wsOverview.Hyperlinks.Add($wsOverview.Range("B5"), "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/1a2b3c4d5e6f7a8b9c0d1e2f3a4b5c6d7e8f9a0b/e2e/d8893c56-fe19-464e-94ef-f52a94eba1ad.md", "", "", $pathAndName) | Out-Null
$wsOverview.Range("B5").Style = "HyperLink"

# -------------------------------------------------------------------
# Sheet 2: zh-cn
# -------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item(2)
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Range("A5").Value = $fileName
$wsZhCn.Range("B5").Value = $ext
$wsZhCn.Range("C5").Value = $statusInSync
$wsZhCn.Range("D5").Value = "e2e"
$wsZhCn.Range("E5").Value = "ht"
$wsZhCn.Range("F5").Value = "'True"
$wsZhCn.Range("G5").Value = $zhXlf
$wsZhCn.Range("H5").Value = $zhHandoffDate
$wsZhCn.Range("H5").NumberFormat = $dateFmt
$wsZhCn.Range("J5").Value = $fileName
$wsZhCn.Range("K5").Value = $zhXlf
$wsZhCn.Range("L5").Value = $zhHandbackDate
$wsZhCn.Range("L5").NumberFormat = $dateFmt
$wsZhCn.Range("O5").Value = "'True"
$wsZhCn.Range("Q5").Value = "'False"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/1a2b3c4d5e6f7a8b9c0d1e2f3a4b5c6d7e8f9a0b/e2e/d8893c56-fe19-464e-94ef-f52a94eba1ad.md", "", "", $fileName) | Out-Null
$wsZhCn.Range("A5").Style = "HyperLink"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("J5"), "https://github.com/OpenLocalizationTestOrg/ol-test4-zhcn/blob/2b3c4d5e6f7a8b9c0d1e2f3a4b5c6d7e8f9a0b1c/e2e/d8893c56-fe19-464e-94ef-f52a94eba1ad.md", "", "", $fileName) | Out-Null
$wsZhCn.Range("J5").Style = "HyperLink"

# -------------------------------------------------------------------
# Sheet 3: de-de
# -------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item(3)
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Range("A5").Value = $fileName
$wsDeDe.Range("B5").Value = $ext
$wsDeDe.Range("C5").Value = $statusInSync
$wsDeDe.Range("D5").Value = "e2e"
$wsDeDe.Range("E5").Value = "ht"
$wsDeDe.Range("F5").Value = "True"
$wsDeDe.Range("G5").Value = $deXlf
$wsDeDe.Range("H5").Value = $deHandoffDate
$wsDeDe.Range("H5").NumberFormat = $dateFmt
$wsDeDe.Range("J5").Value = $fileName
$wsDeDe.Range("K5").Value = $deXlf
$wsDeDe.Range("L5").Value = $deHandbackDate
$wsDeDe.Range("L5").NumberFormat = $dateFmt
$wsDeDe.Range("O5").Value = "True"
$wsDeDe.Range("Q5").Value = "False"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/1a2b3c4d5e6f7a8b9c0d1e2f3a4b5c6d7e8f9a0b/e2e/d8893c56-fe19-464e-94ef-f52a94eba1ad.md", "", "", $fileName) | Out-Null
$wsDeDe.Range("A5").Style = "HyperLink"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("J5"), "https://github.com/OpenLocalizationTestOrg/ol-test4-dede/blob/3c4d5e6f7a8b9c0d1e2f3a4b5c6d7e8f9a0b1c2d/e2e/d8893c56-fe19-464e-94ef-f52a94eba1ad.md", "", "", $fileName) | Out-Null
$wsDeDe.Range("J5").Style = "HyperLink"
